# LDLC Suivi smartphones - append a new price-check column (Q) that mirrors
# the most recent existing column (P): same values/blank pattern, same
# header style, with a fresh timestamp in the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the whole of column P (header + all data/blank rows) into column Q.
# Using Range.Copy (rather than re-typing every value) carries over the
# header's cell style (bold/border/centered) and keeps the blank rows
# blank, exactly mirroring column P's shape.
$ws.Range("P1:P204").Copy($ws.Range("Q1:Q204"))

# The new column's header gets its own (later) timestamp.
$ws.Range("Q1").Value = "2026-01-28 08:17:13"
